# Adds new columns I ("I0") and J ("IF") to the sheet, mirroring the
# existing header style from column H, and fills in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, border, centered) from H1 onto I1:J1
# so the new header cells share the same style as the existing headers.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$rows = @(
    @{Row=2; I=9; J=9}
    @{Row=3; I=9; J=9}
    @{Row=4; I=8; J=8}
    @{Row=5; I=7; J=8}
    @{Row=6; I=7; J=7}
    @{Row=7; I=8; J=8}
    @{Row=8; I=9; J=9}
    @{Row=9; I=8; J=8}
    @{Row=10; I=8; J=8}
    @{Row=11; I=7; J=7}
    @{Row=12; I=8; J=8}
    @{Row=13; I=9; J=9}
    @{Row=14; I=8; J=8}
    @{Row=15; I=7; J=7}
    @{Row=16; I=8; J=8}
    @{Row=17; I=7; J=7}
    @{Row=18; I=7; J=7}
    @{Row=19; I=8; J=8}
    @{Row=20; I=8; J=8}
    @{Row=21; I=7; J=7}
    @{Row=22; I=8; J=8}
    @{Row=23; I=8; J=8}
    @{Row=24; I=8; J=8}
    @{Row=25; I=7; J=7}
    @{Row=26; I=8; J=8}
    @{Row=27; I=9; J=9}
    @{Row=28; I=9; J=9}
    @{Row=29; I=9; J=9}
    @{Row=30; I=9; J=9}
    @{Row=31; I=9; J=9}
    @{Row=32; I=9; J=9}
    @{Row=33; I=8; J=8}
    @{Row=34; I=8; J=9}
    @{Row=35; I=8; J=8}
    @{Row=36; I=8; J=8}
    @{Row=37; I=9; J=9}
    @{Row=38; I=8; J=8}
    @{Row=39; I=9; J=9}
    @{Row=40; I=9; J=9}
    @{Row=41; I=8; J=8}
    @{Row=42; I=9; J=9}
    @{Row=43; I=8; J=8}
    @{Row=44; I=9; J=9}
    @{Row=45; I=9; J=9}
    @{Row=46; I=9; J=9}
    @{Row=47; I=8; J=8}
    @{Row=48; I=7; J=7}
    @{Row=49; I=9; J=9}
    @{Row=50; I=8; J=8}
    @{Row=51; I=9; J=9}
    @{Row=52; I=8; J=8}
    @{Row=53; I=7; J=7}
    @{Row=54; I=9; J=9}
    @{Row=55; I=9; J=9}
    @{Row=56; I=7; J=7}
    @{Row=57; I=8; J=8}
    @{Row=58; I=5; J=5}
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 9).Value = $entry.I
    $ws.Cells.Item($r, 10).Value = $entry.J
}

Write-Output "I0 and IF columns added"
